$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32: Crab Oil
$ws.Range("H32").Value = 111111940
$ws.Range("I32").Value = 200000690
$ws.Range("J32").Value = 998
$ws.Range("K32").Value = 200000690
$ws.Range("L32").Value = 998
$ws.Range("M32").Value = -200000364
$ws.Range("N32").Value = -1650

# Row 98: Enchanted Durium Ink
$ws.Range("H98").Value = 1356.1364
$ws.Range("I98").Value = 985.2778
$ws.Range("J98").Value = 3025
$ws.Range("K98").Value = 985.2778
$ws.Range("L98").Value = 3025
$ws.Range("M98").Value = 512.7222
$ws.Range("N98").Value = -6021

# Row 122: Enchanted High Durium Ink
$ws.Range("H122").Value = 1356.1364
$ws.Range("I122").Value = 985.2778
$ws.Range("J122").Value = 3025
$ws.Range("K122").Value = 2955.8334
$ws.Range("L122").Value = 9075
$ws.Range("M122").Value = -505.8334
$ws.Range("N122").Value = -13975

# Row 132: Growth Formula Lambda
$ws.Range("H132").Value = 3407.2258
$ws.Range("I132").Value = 3027.4666
$ws.Range("J132").Value = 14800
$ws.Range("K132").Value = 9082.399800000001
$ws.Range("L132").Value = 44400
$ws.Range("M132").Value = -6552.399800000001
$ws.Range("N132").Value = -49460

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 4689000
$ws.Range("I137").Value = 2633089.8
$ws.Range("J137").Value = 7693792
$ws.Range("K137").Value = 7899269.399999999
$ws.Range("L137").Value = 23081376
$ws.Range("M137").Value = -7896719.399999999
$ws.Range("N137").Value = -23086476

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 5760.54
$ws.Range("I32").Value = 3515.8901
$ws.Range("J32").Value = 28456.445
$ws.Range("K32").Value = 3515.8901
$ws.Range("L32").Value = 28456.445
$ws.Range("M32").Value = -3228.8901
$ws.Range("N32").Value = -29030.445

# Row 61: Cobalt Ingot
$ws.Range("H61").Value = 1687.1034
$ws.Range("I61").Value = 1659.0769
$ws.Range("J61").Value = 1930
$ws.Range("K61").Value = 1659.0769
$ws.Range("L61").Value = 1930
$ws.Range("M61").Value = -1447.0769
$ws.Range("N61").Value = -2354

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 2388.8333
$ws.Range("I122").Value = 2138.3845
$ws.Range("J122").Value = 3040
$ws.Range("K122").Value = 6415.1535
$ws.Range("L122").Value = 9120
$ws.Range("M122").Value = -3965.1535
$ws.Range("N122").Value = -14020

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 131163.28
$ws.Range("I132").Value = 174212.97
$ws.Range("J132").Value = 6319.2
$ws.Range("K132").Value = 522638.91
$ws.Range("L132").Value = 18957.6
$ws.Range("M132").Value = -520108.91
$ws.Range("N132").Value = -24017.6

# Row 136: Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1687.1034
$ws.Range("I136").Value = 1659.0769
$ws.Range("J136").Value = 1930
$ws.Range("K136").Value = 4977.2307
$ws.Range("L136").Value = 5790
$ws.Range("M136").Value = -2427.2307
$ws.Range("N136").Value = -10890

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steel Nugget
$ws.Range("H94").Value = 851.625
$ws.Range("I94").Value = 1072.6666
$ws.Range("J94").Value = 719
$ws.Range("K94").Value = 1072.6666
$ws.Range("L94").Value = 719
$ws.Range("M94").Value = -621.6666
$ws.Range("N94").Value = -1621

# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 84839.19500000001
$ws.Range("I134").Value = 108345.82
$ws.Range("J134").Value = 2566
$ws.Range("K134").Value = 325037.46
$ws.Range("L134").Value = 7698
$ws.Range("M134").Value = -322502.46
$ws.Range("N134").Value = -12768

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Frantoio Oil
$ws.Range("H107").Value = 646.875
$ws.Range("I107").Value = 319.6
$ws.Range("K107").Value = 958.8000000000001
$ws.Range("M107").Value = 961.1999999999999

# Row 122: Northern Sea Salt
$ws.Range("H122").Value = 48887.61
$ws.Range("I122").Value = 337.58823
$ws.Range("J122").Value = 60678.33
$ws.Range("K122").Value = 3038.29407
$ws.Range("L122").Value = 546104.97
$ws.Range("M122").Value = -588.2940699999999
$ws.Range("N122").Value = -551004.97

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Durium Ingot
$ws.Range("H102").Value = 2757.5
$ws.Range("I102").Value = 2421
$ws.Range("J102").Value = 4440
$ws.Range("K102").Value = 2421
$ws.Range("L102").Value = 4440
$ws.Range("M102").Value = -799
$ws.Range("N102").Value = -7684

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Leather
$ws.Range("H7").Value = 6668572
$ws.Range("I7").Value = 11112878
$ws.Range("J7").Value = 2113.3333
$ws.Range("K7").Value = 11112878
$ws.Range("L7").Value = 2113.3333
$ws.Range("M7").Value = -11112766
$ws.Range("N7").Value = -2337.3333

# Row 16: Hard Leather
$ws.Range("H16").Value = 350.2857
$ws.Range("I16").Value = 350.2857
$ws.Range("K16").Value = 350.2857
$ws.Range("M16").Value = -180.2857

# Row 55: Peiste Leather
$ws.Range("H55").Value = 411.30768
$ws.Range("I55").Value = 228.55556
$ws.Range("J55").Value = 508.05884
$ws.Range("K55").Value = 228.55556
$ws.Range("L55").Value = 508.05884
$ws.Range("M55").Value = -55.55556000000001
$ws.Range("N55").Value = -854.0588399999999

# Row 93: Gagana Leather
$ws.Range("H93").Value = 1411.5333
$ws.Range("I93").Value = 1121.9
$ws.Range("K93").Value = 1121.9
$ws.Range("M93").Value = 126.0999999999999

# Row 100: Tiger Leather
$ws.Range("H100").Value = 1530.5555
$ws.Range("I100").Value = 1379.1666
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 1379.1666
$ws.Range("L100").Value = 1833.3334
$ws.Range("M100").Value = -838.1666
$ws.Range("N100").Value = -2915.3334

# Row 126: Saiga Leather
$ws.Range("H126").Value = 6668572
$ws.Range("I126").Value = 11112878
$ws.Range("J126").Value = 2113.3333
$ws.Range("K126").Value = 33338634
$ws.Range("L126").Value = 6339.999899999999
$ws.Range("M126").Value = -33336164
$ws.Range("N126").Value = -11279.9999

# Row 136: Br'aax Leather
$ws.Range("H136").Value = 1520.2554
$ws.Range("I136").Value = 1411.425
$ws.Range("J136").Value = 2142.1428
$ws.Range("K136").Value = 4234.275
$ws.Range("L136").Value = 6426.428400000001
$ws.Range("M136").Value = -1684.275
$ws.Range("N136").Value = -11526.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Ramie Robe of Casting
$ws.Range("H74").Value = 8417
$ws.Range("J74").Value = 8417
$ws.Range("L74").Value = 8417
$ws.Range("N74").Value = -10289

# Row 77: Ramie Robe of Casting
$ws.Range("H77").Value = 8417
$ws.Range("J77").Value = 8417
$ws.Range("L77").Value = 25251
$ws.Range("N77").Value = -34611

# Row 81: Crawler Silk
$ws.Range("H81").Value = 6897.0527
$ws.Range("I81").Value = 10458.546
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 20917.092
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -19856.092
$ws.Range("N81").Value = -6122

# Row 84: Crawler Silk
$ws.Range("H84").Value = 6897.0527
$ws.Range("I84").Value = 10458.546
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 104585.46
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -99281.46000000001
$ws.Range("N84").Value = -30608

# Row 113: Pixie Floss
$ws.Range("H113").Value = 44535.695
$ws.Range("I113").Value = 67078.07000000001
$ws.Range("J113").Value = 2268.75
$ws.Range("K113").Value = 201234.21
$ws.Range("L113").Value = 6806.25
$ws.Range("M113").Value = -199064.21
$ws.Range("N113").Value = -11146.25

# Row 122: Dark Hempen Cloth
$ws.Range("H122").Value = 10818
$ws.Range("I122").Value = 12121.789
$ws.Range("J122").Value = 4625
$ws.Range("K122").Value = 36365.367
$ws.Range("L122").Value = 13875
$ws.Range("M122").Value = -33915.367
$ws.Range("N122").Value = -18775

# Row 126: Snow Linen
$ws.Range("H126").Value = 1196.8948
$ws.Range("I126").Value = 729.4
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 2188.2
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = 281.8000000000002
$ws.Range("N126").Value = -13790

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 1911.8667
$ws.Range("I132").Value = 2164.8538
$ws.Range("J132").Value = 1365.9474
$ws.Range("K132").Value = 6494.5614
$ws.Range("L132").Value = 4097.8422
$ws.Range("M132").Value = -3964.5614
$ws.Range("N132").Value = -9157.842199999999

# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 1390.7122
$ws.Range("I136").Value = 1357.1177
$ws.Range("J136").Value = 1504.9333
$ws.Range("K136").Value = 4071.3531
$ws.Range("L136").Value = 4514.7999
$ws.Range("M136").Value = -1521.3531
$ws.Range("N136").Value = -9614.7999
